$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.209.84'
$ws.Range('E2').Value = '  +2.57%  '
$ws.Range('D3').Value = '1.915.60'
$ws.Range('E3').Value = '  +2.25%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.006'
$ws.Range('E4').Value = '  -0.93%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '317.36'
$ws.Range('E5').Value = '  +1.30%  '
$ws.Range('E6').Value = '  -0.90%  '
$ws.Range('E7').Value = '  +0.87%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3823'
$ws.Range('E8').Value = '  +1.47%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07381'
$ws.Range('E9').Value = '  -0.01%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.9401'
$ws.Range('E10').Value = '  -0.02%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '20.87'
$ws.Range('E11').Value = '  +1.00%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07829'
$ws.Range('E12').Value = '  -0.57%  '
$ws.Range('D13').Value = '1.941.38'
$ws.Range('E13').Value = '  +3.56%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.506'
$ws.Range('E14').Value = '  +1.01%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.642'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '91.27'
$ws.Range('E16').Value = '  +0.29%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.007'
$ws.Range('E17').Value = '  -0.94%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008840'
$ws.Range('E18').Value = '  -1.50%  '
$ws.Range('D20').Value = '28.227.06'
$ws.Range('E21').Value = '  -0.46%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.164'
$ws.Range('E22').Value = '  +0.37%  '
$ws.Range('D23').Value = '2.152.26'
$ws.Range('E23').Value = '  +1.87%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '10.93'
$ws.Range('E24').Value = '  +1.98%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '156.37'
$ws.Range('E25').Value = '  +1.55%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.925'
$ws.Range('E26').Value = '  -1.24%  '
$ws.Range('E27').Value = '  +0.04%  '
$ws.Range('E28').Value = '  +4.08%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '116.36'
$ws.Range('E29').Value = '  +0.18%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.970'
$ws.Range('E30').Value = '  -0.63%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08903'
$ws.Range('E31').Value = '  -0.29%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.361'
$ws.Range('E32').Value = '  +0.98%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.251'
$ws.Range('E33').Value = '  +2.63%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7738'
$ws.Range('E34').Value = '  +3.05%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.702'
$ws.Range('E35').Value = '  +2.12%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.659'
$ws.Range('E36').Value = '  -1.96%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02053'
$ws.Range('E37').Value = '  -1.07%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.102'
$ws.Range('E38').Value = '  -1.74%  '
$ws.Range('B39').Value = 'Hedera'
$ws.Range('C39').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.05322'
$ws.Range('E39').Value = '  +0.30%  '
$ws.Range('B40').Value = 'TheSandbox'
$ws.Range('C40').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.5526'
$ws.Range('E40').Value = '  +2.77%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.008'
$ws.Range('E42').Value = '  -0.59%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.1530'
$ws.Range('E43').Value = '  +0.37%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.466'
$ws.Range('E44').Value = '  +0.26%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '10.71'
$ws.Range('E45').Value = '  +0.91%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.4862'
$ws.Range('E46').Value = '  +0.31%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '106.94'
$ws.Range('E47').Value = '  +3.65%  '
$ws.Range('E48').Value = '  -0.96%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '68.72'
$ws.Range('E50').Value = '  +2.26%  '
$ws.Range('E51').Value = '  +0.05%  '
